$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.860.04'
$ws.Range('E2').Value = '  +1.53%  '
$ws.Range('D3').Value = '1.763.76'
$ws.Range('E3').Value = '  +1.53%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '327.94'
$ws.Range('E5').Value = '  +1.62%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4463'
$ws.Range('E7').Value = '  -1.85%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3542'
$ws.Range('E8').Value = '  +0.45%  '
$ws.Range('E9').Value = '  +1.64%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.07404'
$ws.Range('E10').Value = '  +0.34%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '1.099'
$ws.Range('E11').Value = '  +2.20%  '
$ws.Range('E12').Value = '  +0.07%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '20.91'
$ws.Range('E13').Value = '  +2.58%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.020'
$ws.Range('E14').Value = '  +1.86%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.230'
$ws.Range('E15').Value = '  +2.63%  '
$ws.Range('D16').Value = '1.762.89'
$ws.Range('E16').Value = '  +1.76%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '92.97'
$ws.Range('E17').Value = '  +2.10%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.00001061'
$ws.Range('E18').Value = '  +0.77%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06425'
$ws.Range('E19').Value = '  +1.46%  '
$ws.Range('E20').Value = '  +0.07%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '17.09'
$ws.Range('E21').Value = '  +3.00%  '
$ws.Range('E22').Value = '  +0.17%  '
$ws.Range('D23').Value = '27.905.18'
$ws.Range('E23').Value = '  +1.55%  '
$ws.Range('E24').Value = '  +0.93%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.107'
$ws.Range('E25').Value = '  +1.66%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '160.99'
$ws.Range('E26').Value = '  -0.32%  '
$ws.Range('E27').Value = '  +1.96%  '
$ws.Range('D28').Value = '1.968.30'
$ws.Range('E28').Value = '  +1.88%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.152'
$ws.Range('E29').Value = '  +5.15%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '124.56'
$ws.Range('E30').Value = '  -0.26%  '
$ws.Range('E31').Value = '  +5.27%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.09189'
$ws.Range('E32').Value = '  +0.90%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '5.638'
$ws.Range('E33').Value = '  +4.59%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.690'
$ws.Range('E34').Value = '  +1.13%  '
$ws.Range('E35').Value = '  +2.13%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.06197'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.02280'
$ws.Range('E37').Value = '  +0.65%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.2099'
$ws.Range('E38').Value = '  +2.32%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.6310'
$ws.Range('E39').Value = '  +1.40%  '
$ws.Range('E40').Value = '  +1.43%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.188'
$ws.Range('E41').Value = '  +0.24%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.394'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '7.862'
$ws.Range('E43').Value = '  +2.34%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '13.22'
$ws.Range('E44').Value = '  +1.53%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.743'
$ws.Range('E45').Value = '  +1.30%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.5856'
$ws.Range('E46').Value = '  +1.22%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '122.38'
$ws.Range('E47').Value = '  +0.45%  '
$ws.Range('E48').Value = '  +1.51%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.06899'
$ws.Range('E49').Value = '  +0.89%  '
$ws.Range('E50').Value = '  +2.05%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '72.85'
$ws.Range('E51').Value = '  +2.58%  '
